$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "1.00", "61.296.35") are preserved verbatim as text, matching
# the source inlineStr cells instead of being coerced to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.296.35'
$ws.Range("E2").Value = '  -4.64%  '
$ws.Range("D3").Value = '3.314.82'
$ws.Range("E3").Value = '  -5.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '566.52'
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").Value = '128.85'
$ws.Range("E6").Value = '  -2.99%  '
$ws.Range("D8").Value = '3.313.83'
$ws.Range("E8").Value = '  -5.13%  '
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").Value = '7.31'
$ws.Range("E10").Value = '  -5.26%  '
$ws.Range("E11").Value = '  -4.16%  '
$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("D13").Value = '3.882.13'
$ws.Range("E13").Value = '  -5.04%  '
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").Value = '3.320.09'
$ws.Range("E15").Value = '  -4.96%  '
$ws.Range("E16").Value = '  -5.48%  '
$ws.Range("D17").Value = '24.66'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").Value = '61.410.86'
$ws.Range("E18").Value = '  -4.37%  '
$ws.Range("D19").Value = '5.68'
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '9.00'
$ws.Range("E21").Value = '  -10.27%  '
$ws.Range("D22").Value = '355.18'
$ws.Range("E22").Value = '  -7.76%  '
$ws.Range("D23").Value = '0.558'
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '3.447.09'
$ws.Range("E25").Value = '  -5.13%  '
$ws.Range("D26").Value = '69.37'
$ws.Range("E26").Value = '  -6.74%  '
$ws.Range("E27").Value = '  -5.51%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '7.22'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").Value = '1.45'
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  -5.82%  '
$ws.Range("E33").Value = '  -2.90%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '3.343.55'
$ws.Range("E35").Value = '  -5.10%  '
$ws.Range("D36").Value = '22.57'
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("D37").Value = '5.32'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '6.82'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").Value = '161.12'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").Value = '0.0759'
$ws.Range("E41").Value = '  -2.82%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '4.39'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").Value = '41.07'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("E45").Value = '  -7.57%  '
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").Value = '1.55'
$ws.Range("E47").Value = '  -4.90%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '22.26'
$ws.Range("E48").Value = '  -7.75%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '6.72'
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").Value = '0.858'
$ws.Range("E50").Value = '  -8.10%  '
$ws.Range("D51").Value = '21.36'
$ws.Range("E51").Value = '  +3.00%  '

# Restore the original (default) style so no stray formatting is
# left behind on the data range.
$ws.Range("B2:E51").Style = "Normal"
